$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 15

$ws.Cells.Item($r, 1).Value = 131082500          # A  Id
$ws.Cells.Item($r, 2).Value = 58043              # B  Taxonsorteringsordning
$ws.Cells.Item($r, 4).Value = "NT"               # D  Rödlistade
$ws.Cells.Item($r, 5).Value = 103021             # E  TaxonId
$ws.Cells.Item($r, 6).Value = "Talltita"         # F  Artnamn
$ws.Cells.Item($r, 7).Value = "Poecile montanus" # G  Vetenskapligt namn
$ws.Cells.Item($r, 8).Value = "(Conrad von Baldenstein, 1827)" # H  Auktor

# I  Antal - stored as text "2" (matches pattern used elsewhere in this column)
$ws.Cells.Item($r, 9).NumberFormat = "@"
$ws.Cells.Item($r, 9).Value = "2"

$ws.Cells.Item($r, 11).Value = "adult"                 # K  Ålder-Stadium
$ws.Cells.Item($r, 13).Value = "permanent revir"       # M  Aktivitet
$ws.Cells.Item($r, 16).Value = "Parsen, naturskogsrest Holmen Skog, Ög" # P Lokalnamn
$ws.Cells.Item($r, 17).Value = 571116            # Q  Ost
$ws.Cells.Item($r, 18).Value = 6467411           # R  Nord
$ws.Cells.Item($r, 19).Value = 10                # S  Noggrannhet
$ws.Cells.Item($r, 20).Value = "Östergötland"    # T  Län
$ws.Cells.Item($r, 21).Value = "Söderköping"     # U  Kommun
$ws.Cells.Item($r, 22).Value = "Östergötland"    # V  Provins
$ws.Cells.Item($r, 23).Value = "Östra Ryd"       # W  Socken

$ws.Cells.Item($r, 25).NumberFormat = "@"        # Y  Startdatum (kept as text)
$ws.Cells.Item($r, 25).Value = "2026-02-07"

$ws.Cells.Item($r, 27).NumberFormat = "@"        # AA Slutdatum (kept as text)
$ws.Cells.Item($r, 27).Value = "2026-02-07"

$ws.Cells.Item($r, 29).Value = "Revirparet"      # AC Publik kommentar
$ws.Cells.Item($r, 30).Value = $false            # AD Ej återfunnen
$ws.Cells.Item($r, 31).Value = $false            # AE Osäker artbestämning
$ws.Cells.Item($r, 33).Value = $false            # AG Ospontan
$ws.Cells.Item($r, 49).Value = "Steve Daurer"    # AW Rapportör
$ws.Cells.Item($r, 50).Value = "Steve Daurer"    # AX Observatörer
